$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet (also updates the defined name '_20160928_UNG' reference) ---
$ws.Name = "20161014 XLF"

# --- Update text/path fields (order matters: it controls the order new
#     shared strings are appended in, which must match the target file) ---
$path = "E:\Datos\bolsa\cuenta personal\analisis de valores\Trades activos\Scanning\20161014"
$ws.Range("B3").Value = $path
$ws.Range("B20").Value = $path
$ws.Range("B2").Value = "20161014 +XLF-161021C21.00"
$ws.Range("B17").Value = "XLF"
$ws.Range("B19").Value = "20161014 +XLF-161021C20.00"

# --- Numeric field updates ---
$ws.Range("B5").Value = 21
$ws.Range("B9").Value = 0.01
$ws.Range("B10").Value = 19.434999999999999
$ws.Range("B12").Value = 10
$ws.Range("B13").Value = 14
$ws.Range("B14").Value = 13
$ws.Range("B15").Value = 33
$ws.Range("B16").Value = 31
$ws.Range("B18").Value = 0.13
$ws.Range("B21").Value = 20
$ws.Range("B22").Value = -0.02
$ws.Range("B23").Value = "19,47,5"
$ws.Range("B26").Value = 14
$ws.Range("B27").Value = 13
$ws.Range("B28").Value = 38
$ws.Range("B29").Value = 27
$ws.Range("B30").Value = 0.13

# --- Apply cell formatting: column B body right-aligned, header centered ---
$ws.Range("B2:B30").HorizontalAlignment = -4152
$ws.Range("B1").HorizontalAlignment = -4108

# --- Update selection to match the saved workbook view ---
$ws.Range("B30").Select()
